$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.08511000000000001
$ws.Range("D2").Value = 0.5594
$ws.Range("E2").Value = 29.39
$ws.Range("F2").Value = 9
$ws.Range("C3").Value = 0.6979
$ws.Range("D3").Value = 0.6743
$ws.Range("E3").Value = 5.193
$ws.Range("F3").Value = 23
$ws.Range("C4").Value = 0.6943
$ws.Range("D4").Value = 0.6717
$ws.Range("E4").Value = 2.554
$ws.Range("F4").Value = 24
$ws.Range("C5").Value = 0.4422
$ws.Range("D5").Value = 1.353
$ws.Range("E5").Value = 6.515
$ws.Range("F5").Value = 19
$ws.Range("C6").Value = 0.3291
$ws.Range("D6").Value = 1.443
$ws.Range("E6").Value = 4.434
$ws.Range("F6").Value = 20
$ws.Range("C7").Value = 0.3268
$ws.Range("D7").Value = 1.498
$ws.Range("E7").Value = 1.808
$ws.Range("F7").Value = 21
$ws.Range("C8").Value = 0.08572
$ws.Range("D8").Value = 1.398
$ws.Range("E8").Value = 12.8
$ws.Range("F8").Value = 12
$ws.Range("C9").Value = 0.2293
$ws.Range("D9").Value = 1.712
$ws.Range("E9").Value = 8.132
$ws.Range("F9").Value = 13
$ws.Range("D10").Value = 1.681
$ws.Range("E10").Value = 8.157
$ws.Range("F10").Value = 18
$ws.Range("C11").Value = 0.5077
$ws.Range("D11").Value = 1.704
$ws.Range("E11").Value = 1.981
$ws.Range("F11").Value = 19
$ws.Range("C12").Value = 0.3098
$ws.Range("D12").Value = 1.638
$ws.Range("E12").Value = 0.2045
$ws.Range("F12").Value = 22
$ws.Range("C13").Value = 0.3686
$ws.Range("D13").Value = 1.656
$ws.Range("E13").Value = 3.84
$ws.Range("F13").Value = 23
$ws.Range("C14").Value = 0.01326
$ws.Range("D14").Value = 2.343
$ws.Range("E14").Value = 1.956
$ws.Range("F14").Value = 12
$ws.Range("C15").Value = 0.1076
$ws.Range("D15").Value = 2.419
$ws.Range("E15").Value = 1.707
$ws.Range("F15").Value = 12
$ws.Range("C16").Value = 0.01744
$ws.Range("D16").Value = 2.7
$ws.Range("E16").Value = 0.7079
$ws.Range("F16").Value = 14
$ws.Range("C17").Value = 0.001435
$ws.Range("D17").Value = 3.088
$ws.Range("E17").Value = 1.282
$ws.Range("F17").Value = 8
$ws.Range("C18").Value = 0.02765
$ws.Range("D18").Value = 2.91
$ws.Range("E18").Value = 1.842
$ws.Range("F18").Value = 12
$ws.Range("C19").Value = 0.2279
$ws.Range("D19").Value = 2.772
$ws.Range("E19").Value = 1.935
$ws.Range("F19").Value = 20
$ws.Range("C20").Value = 0.03506
$ws.Range("D20").Value = 3.269
$ws.Range("E20").Value = 7.147
$ws.Range("F20").Value = 13
$ws.Range("C21").Value = 0.03532
$ws.Range("D21").Value = 3.488
$ws.Range("E21").Value = 0.7301
$ws.Range("F21").Value = 14
$ws.Range("C22").Value = 0.005181
$ws.Range("D22").Value = 3.814
$ws.Range("E22").Value = 4.715
$ws.Range("F22").Value = 11
$ws.Range("C23").Value = 0.262
$ws.Range("D23").Value = 2.95
$ws.Range("E23").Value = 1.506
$ws.Range("C24").Value = 0.1853
$ws.Range("D24").Value = 3.051
$ws.Range("E24").Value = 1.543
$ws.Range("C25").Value = 0.001472
$ws.Range("D25").Value = 4.944
$ws.Range("E25").Value = 3.34
$ws.Range("F25").Value = 7
$ws.Range("C26").Value = 0.01113
$ws.Range("D26").Value = 4.415
$ws.Range("E26").Value = 0.1375
$ws.Range("F26").Value = 13
$ws.Range("C27").Value = 0.006732
$ws.Range("D27").Value = 4.759
$ws.Range("E27").Value = 2.331
$ws.Range("F27").Value = 10
$ws.Range("C28").Value = 0.02502
$ws.Range("D28").Value = 4.724
$ws.Range("E28").Value = 0.4165
$ws.Range("F28").Value = 13
$ws.Range("C29").Value = 0.01393
$ws.Range("D29").Value = 4.784
$ws.Range("E29").Value = 2.841
$ws.Range("F29").Value = 13
$ws.Range("C30").Value = 0.001821
$ws.Range("D30").Value = 5.247
$ws.Range("E30").Value = 1.1
$ws.Range("F30").Value = 7
$ws.Range("C31").Value = 0.1107
$ws.Range("D31").Value = 4.248
$ws.Range("E31").Value = 1.844
$ws.Range("F31").Value = 23
$ws.Range("C32").Value = 0.143
$ws.Range("D32").Value = 4.078
$ws.Range("E32").Value = 1.283
$ws.Range("F32").Value = 25
$ws.Range("C33").Value = 0.1074
$ws.Range("D33").Value = 4.168
$ws.Range("E33").Value = 0.5926
$ws.Range("F33").Value = 25
$ws.Range("C34").Value = 0.001119
$ws.Range("D34").Value = 4.909
$ws.Range("E34").Value = 3.076
$ws.Range("F34").Value = 10
$ws.Range("C35").Value = 0.005295
$ws.Range("D35").Value = 4.934
$ws.Range("E35").Value = 1.274
$ws.Range("F35").Value = 12
$ws.Range("C36").Value = 0.0006306
$ws.Range("D36").Value = 5.988
$ws.Range("E36").Value = 0.1078
$ws.Range("F36").Value = 6
$ws.Range("D37").Value = 5.077
$ws.Range("E37").Value = 3.407
$ws.Range("F37").Value = 18
$ws.Range("C38").Value = 0.1844
$ws.Range("D38").Value = 4.817
$ws.Range("E38").Value = 2.233
$ws.Range("F38").Value = 23
$ws.Range("C39").Value = 0.2017
$ws.Range("D39").Value = 4.739
$ws.Range("E39").Value = 0.552
$ws.Range("F39").Value = 24
$ws.Range("C40").Value = 0.1572
$ws.Range("D40").Value = 4.71
$ws.Range("E40").Value = 0.7847
$ws.Range("F40").Value = 25
$ws.Range("C41").Value = 0.01794
$ws.Range("D41").Value = 5.191
$ws.Range("E41").Value = 2.433
$ws.Range("F41").Value = 12
$ws.Range("C42").Value = 0.01892
$ws.Range("D42").Value = 5.273
$ws.Range("E42").Value = 0.6876
$ws.Range("F42").Value = 13
$ws.Range("C43").Value = 0.1761
$ws.Range("D43").Value = 5.155
$ws.Range("E43").Value = 5.025
$ws.Range("F43").Value = 19
$ws.Range("C44").Value = 0.1483
$ws.Range("D44").Value = 4.914
$ws.Range("E44").Value = 6.448
$ws.Range("F44").Value = 24
$ws.Range("C45").Value = 0.1142
$ws.Range("D45").Value = 4.866
$ws.Range("E45").Value = 4.086
$ws.Range("F45").Value = 25
$ws.Range("C46").Value = 0.1509
$ws.Range("D46").Value = 5.466
$ws.Range("E46").Value = 5.745
$ws.Range("F46").Value = 12
$ws.Range("C47").Value = 0.1562
$ws.Range("D47").Value = 4.916
$ws.Range("E47").Value = 5.119
$ws.Range("F47").Value = 25
$ws.Range("C48").Value = 0.07398
$ws.Range("D48").Value = 5.848
$ws.Range("E48").Value = 4.137
$ws.Range("F48").Value = 9
$ws.Range("C49").Value = 0.2396
$ws.Range("D49").Value = 5.176
$ws.Range("E49").Value = 2.015
$ws.Range("F49").Value = 23
$ws.Range("C50").Value = 0.1907
$ws.Range("D50").Value = 4.994
$ws.Range("E50").Value = 1.27
$ws.Range("F50").Value = 25
$ws.Range("C51").Value = 0.002616
$ws.Range("D51").Value = 6.348
$ws.Range("E51").Value = 5.16
$ws.Range("F51").Value = 7
$ws.Range("C52").Value = 0.0207
$ws.Range("D52").Value = 5.95
$ws.Range("E52").Value = 4.706
$ws.Range("F52").Value = 9
$ws.Range("D53").Value = 5.723
$ws.Range("E53").Value = 3.793
$ws.Range("F53").Value = 15
$ws.Range("C54").Value = 0.8924
$ws.Range("D54").Value = 5.521
$ws.Range("E54").Value = 2.419
$ws.Range("F54").Value = 19
$ws.Range("C55").Value = 0.04006
$ws.Range("D55").Value = 5.256
$ws.Range("E55").Value = 7.908
$ws.Range("F55").Value = 9
$ws.Range("C56").Value = 0.06797
$ws.Range("D56").Value = 5.319
$ws.Range("E56").Value = 5.223
$ws.Range("F56").Value = 12
$ws.Range("C57").Value = 0.008366999999999999
$ws.Range("D57").Value = 5.015
$ws.Range("E57").Value = 8.787000000000001
$ws.Range("F57").Value = 9
$ws.Range("C58").Value = 0.01629
$ws.Range("D58").Value = 4.782
$ws.Range("E58").Value = 11.29
$ws.Range("F58").Value = 8
$ws.Range("C59").Value = 0.3051
$ws.Range("D59").Value = 5.119
$ws.Range("E59").Value = 2.438
$ws.Range("F59").Value = 18
$ws.Range("C60").Value = 0.03417
$ws.Range("D60").Value = 4.528
$ws.Range("E60").Value = 5.903
$ws.Range("F60").Value = 12
$ws.Range("C61").Value = 0.05077
$ws.Range("D61").Value = 4.683
$ws.Range("E61").Value = 0.2774
$ws.Range("F61").Value = 14
$ws.Range("D62").Value = 4.784
$ws.Range("E62").Value = 1.055
$ws.Range("F62").Value = 18
$ws.Range("C63").Value = 0.216
$ws.Range("D63").Value = 4.747
$ws.Range("E63").Value = 2.31
$ws.Range("F63").Value = 20
$ws.Range("C64").Value = 0.009783
$ws.Range("D64").Value = 4.563
$ws.Range("E64").Value = 6.836
$ws.Range("F64").Value = 9
$ws.Range("C65").Value = 0.005106
$ws.Range("D65").Value = 5.132
$ws.Range("E65").Value = 6.556
$ws.Range("F65").Value = 7
$ws.Range("C66").Value = 0.02205
$ws.Range("D66").Value = 4.764
$ws.Range("E66").Value = 5.25
$ws.Range("F66").Value = 13
$ws.Range("C67").Value = 0.01528
$ws.Range("D67").Value = 4.891
$ws.Range("E67").Value = 3.808
$ws.Range("F67").Value = 13
$ws.Range("C68").Value = 0.1135
$ws.Range("D68").Value = 4.952
$ws.Range("E68").Value = 3.493
$ws.Range("F68").Value = 17
$ws.Range("C69").Value = 0.1673
$ws.Range("D69").Value = 4.934
$ws.Range("E69").Value = 1.581
$ws.Range("F69").Value = 20
$ws.Range("C70").Value = 0.02316
$ws.Range("D70").Value = 5.525
$ws.Range("E70").Value = 2.211
$ws.Range("F70").Value = 12
$ws.Range("C71").Value = 0.3763
$ws.Range("D71").Value = 4.82
$ws.Range("E71").Value = 0.2145
$ws.Range("F71").Value = 25
$ws.Range("D72").Value = 5.799
$ws.Range("E72").Value = 0.475
$ws.Range("F72").Value = 14
